$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Styling cleanup: remove the bold/border/centered header style ---
# (matches removal of font[1], border[1], and cellXfs[1] in styles.xml)
$ws.Range("A1:U1").ClearFormats()

# --- A1 header label cleared (was "Unnamed: 0") ---
$ws.Range("A1").Value = ""

# --- Row 3: Revisit count corrections ---
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 26
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 25
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 10
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 24

# --- Row 4: Fixation count corrections ---
$ws.Range("D4").Value = 172
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 71
$ws.Range("I4").Value = 1
$ws.Range("K4").Value = 129
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 13
$ws.Range("Q4").Value = 12
$ws.Range("R4").Value = 14
$ws.Range("S4").Value = 151

# --- Row 5: Dwell time (ms) corrections ---
$ws.Range("D5").Value = 43062.2
$ws.Range("E5").Value = 7057.55
$ws.Range("F5").Value = 18953.42
$ws.Range("I5").Value = 166.82
$ws.Range("K5").Value = 32250.73
$ws.Range("N5").Value = 4838.6
$ws.Range("O5").Value = 3937.83
$ws.Range("Q5").Value = 2986.46
$ws.Range("R5").Value = 3286.85
$ws.Range("S5").Value = 39705.59

# --- Row 6: Dwell time (%) corrections ---
$ws.Range("B6").Value = 3.14
$ws.Range("C6").Value = 2.12
$ws.Range("D6").Value = 25.35
$ws.Range("E6").Value = 4.15
$ws.Range("F6").Value = 11.16
$ws.Range("G6").Value = 2.22
$ws.Range("I6").Value = 0.1
$ws.Range("J6").Value = 0.79
$ws.Range("K6").Value = 18.97
$ws.Range("L6").Value = 5.23
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 2.85
$ws.Range("O6").Value = 2.32
$ws.Range("P6").Value = 0.09
$ws.Range("Q6").Value = 1.76
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 23.35
$ws.Range("T6").Value = 0.42
$ws.Range("U6").Value = 3.01

# --- Row 7: Fixation duration (ms) corrections ---
$ws.Range("D7").Value = 250.36
$ws.Range("E7").Value = 252.06
$ws.Range("F7").Value = 266.95
$ws.Range("I7").Value = 166.82
$ws.Range("K7").Value = 250.01
$ws.Range("N7").Value = 403.22
$ws.Range("O7").Value = 302.91
$ws.Range("Q7").Value = 248.87
$ws.Range("R7").Value = 234.77
$ws.Range("S7").Value = 262.95

# --- Row 8: new "First fixation duration (ms)" row populated ---
$ws.Range("A8").Value = "First fixation duration (ms)"
$ws.Range("B8").Value = 83.43
$ws.Range("C8").Value = 333.69
$ws.Range("D8").Value = 183.6
$ws.Range("E8").Value = 433.7
$ws.Range("F8").Value = 100.14
$ws.Range("G8").Value = 333.69
$ws.Range("H8").Value = 69.29
$ws.Range("I8").Value = 166.82
$ws.Range("J8").Value = 150.12
$ws.Range("K8").Value = 100.14
$ws.Range("L8").Value = 417.13
$ws.Range("M8").Value = 166.92
$ws.Range("N8").Value = 300.28
$ws.Range("O8").Value = 100.14
$ws.Range("P8").Value = 150.17
$ws.Range("Q8").Value = 433.7
$ws.Range("R8").Value = 433.78
$ws.Range("S8").Value = 266.96
$ws.Range("T8").Value = 367.05
$ws.Range("U8").Value = 383.73
